$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.914.05"
$ws.Range("E2").Value = "  +2.92%  "

$ws.Range("D3").Value = "3.807.29"
$ws.Range("E3").Value = "  +7.13%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'427.83"
$ws.Range("E5").Value = "  +8.25%  "

$ws.Range("D6").Value = "'130.39"
$ws.Range("E6").Value = "  +3.64%  "

$ws.Range("D7").Value = "3.807.71"
$ws.Range("E7").Value = "  +7.54%  "

$ws.Range("D8").Value = "'0.612"
$ws.Range("E8").Value = "  +3.19%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "'0.732"
$ws.Range("E10").Value = "  +5.53%  "

$ws.Range("D11").Value = "'0.158"
$ws.Range("E11").Value = "  +3.19%  "

$ws.Range("D12").Value = "'0.0000335"
$ws.Range("E12").Value = "  -3.91%  "

$ws.Range("D13").Value = "'41.51"
$ws.Range("E13").Value = "  +5.37%  "

$ws.Range("D14").Value = "'10.49"
$ws.Range("E14").Value = "  +12.42%  "

$ws.Range("D15").Value = "4.423.76"
$ws.Range("E15").Value = "  +7.88%  "

$ws.Range("D16").Value = "'15.38"
$ws.Range("E16").Value = "  +20.14%  "

$ws.Range("E17").Value = "  +1.30%  "

$ws.Range("D18").Value = "3.820.35"
$ws.Range("E18").Value = "  +8.04%  "

$ws.Range("D19").Value = "'20.05"
$ws.Range("E19").Value = "  +6.00%  "

$ws.Range("E20").Value = "  +8.19%  "

$ws.Range("D21").Value = "66.171.58"
$ws.Range("E21").Value = "  +3.26%  "

$ws.Range("D22").Value = "'414.96"
$ws.Range("E22").Value = "  +3.31%  "

$ws.Range("D23").Value = "'15.22"
$ws.Range("E23").Value = "  +8.58%  "

$ws.Range("D24").Value = "'85.44"
$ws.Range("E24").Value = "  +4.17%  "

$ws.Range("E25").Value = "  +8.12%  "

$ws.Range("D26").Value = "'36.97"
$ws.Range("E26").Value = "  +7.54%  "

$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "  +12.65%  "

$ws.Range("D28").Value = "'3.29"
$ws.Range("E28").Value = "  +9.15%  "

$ws.Range("D29").Value = "'5.40"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("D30").Value = "'9.29"
$ws.Range("E30").Value = "  +35.18%  "

$ws.Range("E31").Value = "  +18.51%  "

$ws.Range("D32").Value = "'710.32"
$ws.Range("E32").Value = "  +4.57%  "

$ws.Range("E33").Value = "  +12.91%  "

$ws.Range("E34").Value = "  +4.75%  "

$ws.Range("D35").Value = "'5.84"
$ws.Range("E35").Value = "  +40.03%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").Value = "'38.85"
$ws.Range("E37").Value = "  +3.69%  "

$ws.Range("D38").Value = "'0.149"
$ws.Range("E38").Value = "  -1.91%  "

$ws.Range("D39").Value = "'55.81"
$ws.Range("E39").Value = "  +2.92%  "

$ws.Range("D40").Value = "'0.0470"
$ws.Range("E40").Value = "  +6.28%  "

$ws.Range("D41").Value = "0.0₃0722"
$ws.Range("E41").Value = "  +16.25%  "

$ws.Range("E42").Value = "  +2.83%  "

$ws.Range("E43").Value = "  +0.55%  "

$ws.Range("E44").Value = "  +4.21%  "

$ws.Range("D45").Value = "'3.23"
$ws.Range("E45").Value = "  +6.00%  "

$ws.Range("D46").Value = "'3.39"
$ws.Range("E46").Value = "  +9.15%  "

$ws.Range("D47").Value = "'0.323"
$ws.Range("E47").Value = "  +16.07%  "

$ws.Range("D48").Value = "'2.43"
$ws.Range("E48").Value = "  +41.02%  "

$ws.Range("E49").Value = "  +5.59%  "

$ws.Range("D50").Value = "'2.05"
$ws.Range("E50").Value = "  +4.47%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'141.49"
$ws.Range("E51").Value = "  +0.14%  "
